# Daily attendance processing - 2026-01-03 11:30:58
# Rotate the "Recorded By" (column G) comma-separated list of names/emails
# one position to the left: the first entry is moved to the end of the list.
# Rows whose G value has only a single entry (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $cell.Value = [string]::Join(", ", $rotated)
        }
    }
}
